$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (weekly update), shifting existing rows 5-29 down to 6-30
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the week's new record.
# Most columns repeat the values that were already present for this series (Agricola del
# Norte / Arica y Parinacota / Ciruela / Angeleno / Segunda); only the date and the
# price/volume/unit/origin columns change for this week's entry.
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(5, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(5, 4).Value = "2023-03-02"
$ws.Cells.Item(5, 5).Value = 15
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100103
$ws.Cells.Item(5, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(5, 9).Value = 100103002
$ws.Cells.Item(5, 10).Value = "Ciruela"
$ws.Cells.Item(5, 11).Value = "Angeleno"
$ws.Cells.Item(5, 12).Value = "Segunda"
$ws.Cells.Item(5, 13).Value = 400
$ws.Cells.Item(5, 14).Value = 5000
$ws.Cells.Item(5, 15).Value = 6000
$ws.Cells.Item(5, 16).Value = 5750
$ws.Cells.Item(5, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(5, 18).Value = "Región Metropolitana"
$ws.Cells.Item(5, 19).Value = 575
$ws.Cells.Item(5, 20).Value = 10
